$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) from 45189 to 45190 for all data rows (2-199)
for ($r = 2; $r -le 199; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
